$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.885.29'
$ws.Range("E2").Value = '  +0.27%  '
$ws.Range("D3").Value = '1.638.88'
$ws.Range("E3").Value = '  -0.19%  '
$ws.Range("E4").Value = '  -0.51%  '
$ws.Range("D5").Value = '216.96'
$ws.Range("E5").Value = '  -0.72%  '
$ws.Range("D6").Value = '0.511'
$ws.Range("E6").Value = '  +1.84%  '
$ws.Range("E7").Value = '  -0.51%  '
$ws.Range("E9").Value = '  +0.32%  '
$ws.Range("D10").Value = '19.88'
$ws.Range("E10").Value = '  +3.47%  '
$ws.Range("D11").Value = '0.0848'
$ws.Range("E11").Value = '  +0.01%  '
$ws.Range("D12").Value = '1.868.61'
$ws.Range("E12").Value = '  -0.17%  '
$ws.Range("D13").Value = '1.638.41'
$ws.Range("E13").Value = '  -0.11%  '
$ws.Range("E14").Value = '  -0.85%  '
$ws.Range("D15").Value = '0.530'
$ws.Range("E15").Value = '  +0.87%  '
$ws.Range("D16").Value = '67.16'
$ws.Range("E16").Value = '  +3.06%  '
$ws.Range("D17").Value = '26.881.03'
$ws.Range("E17").Value = '  +0.20%  '
$ws.Range("E18").Value = '  -0.55%  '
$ws.Range("D19").Value = '219.44'
$ws.Range("E19").Value = '  +1.62%  '
$ws.Range("E20").Value = '  -0.56%  '
$ws.Range("D21").Value = '6.84'
$ws.Range("E21").Value = '  +3.11%  '
$ws.Range("D22").Value = '4.39'
$ws.Range("E22").Value = '  +0.52%  '
$ws.Range("E23").Value = '  +3.94%  '
$ws.Range("E24").Value = '  -0.25%  '
$ws.Range("E25").Value = '  -0.46%  '
$ws.Range("E26").Value = '  -0.49%  '
$ws.Range("E27").Value = '  +3.24%  '
$ws.Range("E28").Value = '  +0.57%  '
$ws.Range("E29").Value = '  +0.17%  '
$ws.Range("E30").Value = '  -0.85%  '
$ws.Range("E31").Value = '  -0.90%  '
$ws.Range("D32").Value = '3.33'
$ws.Range("E32").Value = '  -1.41%  '
$ws.Range("E33").Value = '  +0.61%  '
$ws.Range("E34").Value = '  +1.09%  '
$ws.Range("D35").Value = '1.268.66'
$ws.Range("E35").Value = '  +0.22%  '
$ws.Range("D36").Value = '2.44'
$ws.Range("E36").Value = '  -0.12%  '
$ws.Range("E37").Value = '  +2.06%  '
$ws.Range("D38").Value = '0.534'
$ws.Range("E38").Value = '  +0.46%  '
$ws.Range("E39").Value = '  +2.05%  '
$ws.Range("E40").Value = '  -0.54%  '
$ws.Range("E41").Value = '  +0.85%  '
$ws.Range("D42").Value = '5.39'
$ws.Range("E42").Value = '  +0.86%  '
$ws.Range("D43").Value = '1.779.09'
$ws.Range("D44").Value = '2.11'
$ws.Range("E44").Value = '  -1.51%  '
$ws.Range("D45").Value = '61.78'
$ws.Range("E45").Value = '  +0.88%  '
$ws.Range("D46").Value = '91.85'
$ws.Range("E46").Value = '  -0.98%  '
$ws.Range("E47").Value = '  -0.97%  '
$ws.Range("E48").Value = '  +3.30%  '
$ws.Range("E49").Value = '  -0.29%  '
$ws.Range("D50").Value = '7.61'
$ws.Range("E50").Value = '  +0.46%  '
$ws.Range("E51").Value = '  -0.28%  '
